# Fix #70: Remove documentation for regex-unknown from rxriskv and add index "sum-all".
#
# This workbook edit:
#  1. Adds a new column F ("sum_all") to the rxriskv sheet, with a value of 1
#     for every data row (2-47).
#  2. Switches the active/selected tab from "hip_ae_hailer" to "rxriskv", and
#     updates each sheet's remembered selection accordingly.

$wb = $excel.ActiveWorkbook

# --- rxriskv: add the "sum_all" index column ---------------------------------
$wsRx = $wb.Worksheets.Item("rxriskv")

$wsRx.Range("F1").Value = "sum_all"

for ($row = 2; $row -le 47; $row++) {
    $wsRx.Cells.Item($row, 6).Value = 1
}

# --- hip_ae_hailer: was the selected tab before this edit ---------------------
$wsHip = $wb.Worksheets.Item("hip_ae_hailer")
$wsHip.Activate() | Out-Null
$wsHip.Range("C15").Select() | Out-Null

# --- rxriskv: becomes the selected/active tab after this edit -----------------
$wsRx.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$wsRx.Range("B17").Select() | Out-Null
